# Apply the authors's edits to the ChequingAccount test-plan workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (C3, merged C3:D3)
$ws.Range("C3").Value = "Nishant Malhotra"

# Row 7 - __init__ / Attributes are set to input values
$ws.Range("F7").Value = "12345, 1010, 500.0, 2023-01-01, -150.0, 0.05`t"
$ws.Range("G7").Value = "Attributes set: 12345, 1010, 500.0, 2023-01-01, -150.0, 0.05"

# Row 8 - __init__ / overdraft limit has invalid type.
$ws.Range("F8").Value = "limit=`"invalid`"`t"
$ws.Range("G8").Value = "Overdraft Limit defaults to -100.0"

# Row 9 - __init__ / overdraft rate has invalid type.
$ws.Range("F9").Value = "rate=`"invalid`"`tOverdraft Rate defaults to 0.05"
$ws.Range("G9").Value = "Overdraft Rate defaults to 0.05"

# Row 10 - __init__ / date created has invalid type
$ws.Range("F10").Value = "date=`"2023-01-01`" (string)`tDate defaults to date.today()"
$ws.Range("G10").Value = "Date defaults to date.today()"

# Row 11 - get_service_charges / balance greater than overdraft limit
$ws.Range("F11").Value = "balance=50.0, limit=-100.0`t"
$ws.Range("G11").Value = "Service Charge = `$0.50 (Base)"

# Row 12 - get_service_charges / balance less than overdraft limit
$ws.Range("F12").Value = "balance=-600.0, limit=-100.0, rate=0.05`t"
$ws.Range("G12").Value = "Service Charge = `$25.50 (Base + Overdraft Fee)"

# Row 13 - get_service_charges / balance equal to overdraft limit
$ws.Range("F13").Value = "balance=-100.0, limit=-100.0`t"
$ws.Range("G13").Value = "Service Charge = `$0.50 (Base)"

# Row 14 - __str__ / appropriate value returned based on attribute values.
$ws.Range("F14").Value = "1234567, 101, 1559.49, 2023-01-01, -15.0, 0.05`t"
$ws.Range("G14").Value = '"Account Number: 1234567 Balance: $1,559.49\nOverdraft Limit: $-15.00 Overdraft Rate: 5.00% Account Type: Chequing"'

# Selection moved to G7 (single cell) in the saved view
$ws.Range("G7").Select()
